# Update Week 17 target depth data for the Colts (Home row) on both the
# OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet (row 2 = "H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 238
$wsOff.Range("C2").Value = 157
$wsOff.Range("D2").Value = 54
$wsOff.Range("E2").Value = 19

# --- DEF sheet (row 2 = "H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 278
$wsDef.Range("C2").Value = 197
$wsDef.Range("D2").Value = 55
$wsDef.Range("E2").Value = 27
$wsDef.Range("F2").Value = 8
$wsDef.Range("G2").Value = 2

$wb.Save()
